$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")

# Add the new "price" column (B) to Лист1
$ws1.Range("B1").Value = "87 900 ₽"
$ws1.Range("B2").Value = "Не найдено"
$ws1.Range("B3").Value = "105 678 ₽"

# Add a new sheet "BonpetData" right after Лист1, mirroring the same data
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "BonpetData"

$newSheet.Range("A1").Value = "АВР-Б-100-2-1"
$newSheet.Range("B1").Value = "87 900 ₽"
$newSheet.Range("A2").Value = "bababab"
$newSheet.Range("B2").Value = "Не найдено"

# "267515" looks numeric, force it to stay text (matches original Лист1!A3)
$newSheet.Range("A3").Value = "'267515"
$newSheet.Range("A3").ClearFormats()
$newSheet.Range("B3").Value = "105 678 ₽"

# Match default Excel page margins used by the authoring tool for this sheet
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Restore Лист1 as the active/selected sheet (tab) like before the edit
$ws1.Activate()
$ws1.Range("A3").Select() | Out-Null
